$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Start Date"), shifting columns E:V to F:W.
# The new column becomes "Nominated Pharmacy Type" with values "P1" for the data rows,
# copying the format of column D (Nominated Pharmacy) to its left.
$dColumnWidth = $ws.Columns("D").ColumnWidth

$ws.Columns("E").Insert()

# Match formatting of the neighbouring "Nominated Pharmacy" column (D) first
$ws.Columns("E").ColumnWidth = $dColumnWidth
$ws.Range("E1:E3").NumberFormat = "@"

$ws.Range("E1").Value = "Nominated Pharmacy Type"
$ws.Range("E2").Value = "P1"
$ws.Range("E3").Value = "P1"

# Refresh the AutoFilter so its range grows to include the new column
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:W3").AutoFilter()

# The hidden _FilterDatabase defined name keeps the stale range; fix it up
$filterDbName = $ws.Name + "!_FilterDatabase"
$newRef = "=" + $ws.Name + "!`$A`$1:`$W`$3"
foreach ($n in $wb.Names) {
    if ($n.Name -eq $filterDbName) {
        $n.RefersTo = $newRef
    }
}
